# update weights and edu
# The "I:Education" row (row 14) had its weight/covariate name changed
# from "hv108_cont_scale" to "hv106_fctb".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "hv106_fctb"

# Move the active selection to B14 (matches the saved sheetView selection).
$ws.Range("B14").Select()
